$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).Value = 44769
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino dulce"
$ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 60
$ws.Cells.Item($row, 11).Value = 17000
$ws.Cells.Item($row, 12).Value = 18000
$ws.Cells.Item($row, 13).Value = 17500
$ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 972
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"
